$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Z41") holds a TRUE/empty flag. Blank cells were being read by
# pandas as NaN, which upgraded the whole column to float64 and turned the
# TRUE values into 1.0. Give the blank cells an explicit boolean FALSE (with
# a TRUE/FALSE custom display format) so the column stays boolean end to end.
$boolFormat = """TRUE"";""TRUE"";""FALSE"""

$blankFlagRows = @(2, 4, 5, 6, 7, 11, 12)

foreach ($r in $blankFlagRows) {
    $cell = $ws.Range("F$r")
    $cell.NumberFormat = $boolFormat
    $cell.Value = $false
}

# Move the cursor to where the editor left off.
[void]$ws.Range("E29").Select()
